$wb = $excel.ActiveWorkbook

# Insert a new "Token" worksheet right after "Comment" (so the sheet order
# becomes: Comment, Token, Date, Number, Formula, Umlaute).
$commentSheet = $wb.Worksheets.Item("Comment")
$tokenSheet = $wb.Worksheets.Add($null, $commentSheet)
$tokenSheet.Name = "Token"

# Populate the new sheet with its markup-token sample data.
$tokenSheet.Range("A2").Value = "regular"
$tokenSheet.Range("A3").Value = "table"
$tokenSheet.Range("B3").Value = "regular"
$tokenSheet.Range("A5").Value = "end"

# "table" gets a thin box border around it, like the header cells on the
# other sheets.
$tokenSheet.Range("A3").Borders.LineStyle = 1

# Leave the cursor parked on A6, matching the saved selection state.
[void]$tokenSheet.Range("A6").Select()
